$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append new row 15 with the latest mail log entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A15").Value = "Retour status"
$logs.Range("B15").Value = "mailmind.test@zohomail.eu"
$logs.Range("D15").Value = "Retour / Terugbetaling"
$logs.Range("F15").Value = "2025-08-28 20:35:36"
$logs.Range("G15").Value = "Ja"
$logs.Range("H15").Value = "Nee"
$logs.Range("I15").Value = "Nee"
$logs.Range("J15").Value = "Nee"

# --- Extend the conditional-formatting ranges to include the new row ---
$logs.Range("D2:D14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D15"))
$logs.Range("G2:G14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G15"))
$logs.Range("H2:H14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H15"))
$logs.Range("I2:I14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I15"))
$logs.Range("J2:J14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J15"))

# --- "Dashboard" sheet: bump the tally for "Retour / Terugbetaling" ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 14
